$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.268183708190918
$ws.Range("B1").Value = 2.260169982910156
$ws.Range("C1").Value = 6.250511169433594
$ws.Range("D1").Value = 1.496893048286438
$ws.Range("E1").Value = 1.355165123939514
